# "adding 2 more test cases"
# Append two new rows (TestCase_A14, TestCase_A15) to the "Test Cases" sheet,
# mirroring the formatting of the existing rows, then move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Clone the formatting (borders/fill/etc.) of the last existing data row down
# onto the two new rows before filling in their values.
$ws.Range("A9:D9").Copy() | Out-Null
$ws.Range("A15:D15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A15").Value = "TestCase_A14"
$ws.Range("B15").Value = "To verify that user is not able to submit new TR user registration form without filling in the required fields"
$ws.Range("C15").Value = "Y"

$ws.Range("A16").Value = "TestCase_A15"
$ws.Range("B16").Value = "To verify that app doesn't allow the user to create a new account with an email id that has already been used"
$ws.Range("C16").Value = "Y"

# Move the active selection, as recorded in the saved view state.
$ws.Range("A18").Select() | Out-Null
